# Add two new products ("World Match" and "Inplay Matrix") to the
# ProductList / CAPTSetting / SMAPTSetting sheets, mirroring the existing
# rows 2-16 pattern (ProdID 1..15 -> 16 and 17), then leave the workbook's
# selection/active-tab state pointing at the newly added rows, matching
# where the author ended up working (SMAPTSetting, C17:E18 selected).

$wb = $excel.ActiveWorkbook

# --- ProductList: two new product names ---
$ws2 = $wb.Worksheets.Item("ProductList")
$ws2.Range("A17").Value = 16
$ws2.Range("A17").NumberFormat = "0"
$ws2.Range("B17").Value = "World Match"
$ws2.Range("A18").Value = 17
$ws2.Range("A18").NumberFormat = "0"
$ws2.Range("B18").Value = "Inplay Matrix"

# --- CAPTSetting: matching settings rows for the two new products ---
$ws3 = $wb.Worksheets.Item("CAPTSetting")
$ws3.Range("A17").Value = 16
$ws3.Range("A17").NumberFormat = "0"
$ws3.Range("B17").Value = "World Match"
$ws3.Range("C17").Value = 50.5
$ws3.Range("D17").Value = 49
$ws3.Range("E17").Value = 100
$ws3.Range("A18").Value = 17
$ws3.Range("A18").NumberFormat = "0"
$ws3.Range("B18").Value = "Inplay Matrix"
$ws3.Range("C18").Value = 50.5
$ws3.Range("D18").Value = 49
$ws3.Range("E18").Value = 100

# --- SMAPTSetting: matching settings rows for the two new products ---
$ws4 = $wb.Worksheets.Item("SMAPTSetting")
$ws4.Range("A17").Value = 16
$ws4.Range("A17").NumberFormat = "0"
$ws4.Range("B17").Value = "World Match"
$ws4.Range("C17").Value = 20
$ws4.Range("D17").Value = 19
$ws4.Range("E17").Value = 80
$ws4.Range("A18").Value = 17
$ws4.Range("A18").NumberFormat = "0"
$ws4.Range("B18").Value = "Inplay Matrix"
$ws4.Range("C18").Value = 20
$ws4.Range("D18").Value = 19
$ws4.Range("E18").Value = 80

# --- Leave each sheet's selection on the rows that were just added ---
[void]$ws2.Activate()
[void]$ws2.Range("A17:B18").Select()

[void]$ws3.Activate()
[void]$ws3.Range("A17:B18").Select()

# SMAPTSetting is activated last, so it becomes the workbook's active tab.
[void]$ws4.Activate()
[void]$ws4.Range("C17:E18").Select()
